$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Email is already registered."
$ws.Range("C23").Value = 1

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Email or password is incorrect."
$ws.Range("C24").Value = 1

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Registration Successful."
$ws.Range("C25").Value = 1

$ws.Range("E24").Select()
